$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "49.476.99"
$dCell.Style = $dStyle
$ws.Range("E2").Value = "  -0.88%  "

$dCell = $ws.Range("D3")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.634.46"
$dCell.Style = $dStyle
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("E4").Value = "  +0.07%  "

$dCell = $ws.Range("D5")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "112.21"
$dCell.Style = $dStyle
$ws.Range("E5").Value = "  -1.58%  "

$dCell = $ws.Range("D6")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "325.93"
$dCell.Style = $dStyle
$ws.Range("E6").Value = "  -0.40%  "

$dCell = $ws.Range("D7")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.523"
$dCell.Style = $dStyle
$ws.Range("E7").Value = "  -1.29%  "

$dCell = $ws.Range("D8")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.Style = $dStyle
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -0.93%  "

$dCell = $ws.Range("D10")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "39.55"
$dCell.Style = $dStyle
$ws.Range("E10").Value = "  -3.68%  "

$dCell = $ws.Range("D11")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "19.94"
$dCell.Style = $dStyle
$ws.Range("E11").Value = "  -0.99%  "

$dCell = $ws.Range("D12")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.0813"
$dCell.Style = $dStyle
$ws.Range("E12").Value = "  -1.02%  "

$ws.Range("E13").Value = "  +1.67%  "

$dCell = $ws.Range("D14")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "7.63"
$dCell.Style = $dStyle
$ws.Range("E14").Value = "  +3.48%  "

$dCell = $ws.Range("D15")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "3.041.57"
$dCell.Style = $dStyle
$ws.Range("E15").Value = "  -0.56%  "

$dCell = $ws.Range("D16")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.629.84"
$dCell.Style = $dStyle
$ws.Range("E16").Value = "  -1.23%  "

$dCell = $ws.Range("D17")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.854"
$dCell.Style = $dStyle
$ws.Range("E17").Value = "  -2.10%  "

$dCell = $ws.Range("D18")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "49.423.18"
$dCell.Style = $dStyle
$ws.Range("E18").Value = "  -0.80%  "

$dCell = $ws.Range("D19")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "13.39"
$dCell.Style = $dStyle
$ws.Range("E19").Value = "  +1.98%  "

$dCell = $ws.Range("D20")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.91"
$dCell.Style = $dStyle
$ws.Range("E20").Value = "  -0.73%  "

$dCell = $ws.Range("D21")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "6.66"
$dCell.Style = $dStyle
$ws.Range("E21").Value = "  -1.91%  "

$dCell = $ws.Range("D22")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0948"
$dCell.Style = $dStyle
$ws.Range("E22").Value = "  -0.95%  "

$dCell = $ws.Range("D23")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "268.59"
$dCell.Style = $dStyle
$ws.Range("E23").Value = "  -3.24%  "

$dCell = $ws.Range("D24")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "69.10"
$dCell.Style = $dStyle
$ws.Range("E24").Value = "  -4.09%  "

$dCell = $ws.Range("D25")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.56"
$dCell.Style = $dStyle
$ws.Range("E25").Value = "  -1.43%  "

$ws.Range("E26").Value = "  +0.05%  "

$dCell = $ws.Range("D27")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "26.00"
$dCell.Style = $dStyle
$ws.Range("E27").Value = "  -2.89%  "

$dCell = $ws.Range("D28")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "10.13"
$dCell.Style = $dStyle
$ws.Range("E28").Value = "  +1.47%  "

$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("E30").Value = "  -3.02%  "

$dCell = $ws.Range("D31")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "34.47"
$dCell.Style = $dStyle
$ws.Range("E31").Value = "  -4.73%  "

$dCell = $ws.Range("D32")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "49.59"
$dCell.Style = $dStyle
$ws.Range("E32").Value = "  -1.21%  "

$dCell = $ws.Range("D33")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "5.47"
$dCell.Style = $dStyle
$ws.Range("E33").Value = "  +0.57%  "

$dCell = $ws.Range("D34")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.0812"
$dCell.Style = $dStyle
$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("E35").Value = "  -0.16%  "

$dCell = $ws.Range("D36")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "19.06"
$dCell.Style = $dStyle
$ws.Range("E36").Value = "  -2.45%  "

$dCell = $ws.Range("D37")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "4.95"
$dCell.Style = $dStyle
$ws.Range("E37").Value = "  +2.18%  "

$dCell = $ws.Range("D38")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.03"
$dCell.Style = $dStyle
$ws.Range("E38").Value = "  -2.19%  "

$dCell = $ws.Range("D39")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "3.10"
$dCell.Style = $dStyle
$ws.Range("E39").Value = "  +0.09%  "

$dCell = $ws.Range("D40")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "129.52"
$dCell.Style = $dStyle
$ws.Range("E40").Value = "  +3.72%  "

$dCell = $ws.Range("D41")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.37"
$dCell.Style = $dStyle
$ws.Range("E41").Value = "  +6.11%  "

$dCell = $ws.Range("D42")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "22.92"
$dCell.Style = $dStyle
$ws.Range("E42").Value = "  +3.70%  "

$ws.Range("E43").Value = "  -0.99%  "

$dCell = $ws.Range("D44")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.0336"
$dCell.Style = $dStyle
$ws.Range("E44").Value = "  +6.67%  "

$dCell = $ws.Range("D45")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.064.73"
$dCell.Style = $dStyle
$ws.Range("E45").Value = "  -0.65%  "

$dCell = $ws.Range("D46")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "3.30"
$dCell.Style = $dStyle
$ws.Range("E46").Value = "  -0.73%  "

$dCell = $ws.Range("D47")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.12"
$dCell.Style = $dStyle
$ws.Range("E47").Value = "  +7.81%  "

$ws.Range("E48").Value = "  -7.51%  "

$dCell = $ws.Range("D49")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "8.87"
$dCell.Style = $dStyle
$ws.Range("E49").Value = "  -2.83%  "

$dCell = $ws.Range("D50")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "5.21"
$dCell.Style = $dStyle
$ws.Range("E50").Value = "  -3.23%  "

$dCell = $ws.Range("D51")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "58.56"
$dCell.Style = $dStyle
$ws.Range("E51").Value = "  -1.86%  "

